# Generate Report for Handback
# ------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    on the Overview sheet (B/C columns) and on each language sheet's
#    Status column (C).
# 2) Stamp "Latest Handback DateTime" (column H) for both rows on the
#    zh-cn and de-de sheets.
# 3) Populate "Latest Target File" (F) and "Latest Handback File" (G)
#    for both rows on the zh-cn and de-de sheets, each as a hyperlink
#    mirroring the matching Source (A) / Latest Handoff File (D) link.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# Latest Handback DateTime stamps
$wsZhCn.Range("H2").Value = "2016-03-24 18:16:29"
$wsZhCn.Range("H3").Value = "2016-03-24 18:16:29"

$wsDeDe.Range("H2").Value = "2016-03-24 18:16:38"
$wsDeDe.Range("H3").Value = "2016-03-24 18:16:38"

# ------------------------------------------------------------------
# zh-cn sheet: Latest Target File (F) / Latest Handback File (G)
# ------------------------------------------------------------------
$wsZhCn.Range("F2").Value = "23d4fd98-57d0-4775-a293-fe495abf9569.md"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/656fb5ebc286263f4ec67e3f3953099ed9a36ebf/e2e/23d4fd98-57d0-4775-a293-fe495abf9569.md",
    [Type]::Missing,
    [Type]::Missing,
    "23d4fd98-57d0-4775-a293-fe495abf9569.md"
) | Out-Null

$wsZhCn.Range("G2").Value = "23d4fd98-57d0-4775-a293-fe495abf9569.891c8a52f5279b6d8910c11a3544eaa27790e0f8.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4fbded06662da38b74c00b9278aa1eec96323457/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/23d4fd98-57d0-4775-a293-fe495abf9569.891c8a52f5279b6d8910c11a3544eaa27790e0f8.zh-cn.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "23d4fd98-57d0-4775-a293-fe495abf9569.891c8a52f5279b6d8910c11a3544eaa27790e0f8.zh-cn.xlf"
) | Out-Null

$wsZhCn.Range("F3").Value = "c28f4a4c-2e03-4e0a-b9c5-bb6b51823999.md"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/656fb5ebc286263f4ec67e3f3953099ed9a36ebf/e2e/c28f4a4c-2e03-4e0a-b9c5-bb6b51823999.md",
    [Type]::Missing,
    [Type]::Missing,
    "c28f4a4c-2e03-4e0a-b9c5-bb6b51823999.md"
) | Out-Null

$wsZhCn.Range("G3").Value = "c28f4a4c-2e03-4e0a-b9c5-bb6b51823999.6a86a0befb441ebc87fb0e036ebc5d5587d6ddb4.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4fbded06662da38b74c00b9278aa1eec96323457/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/c28f4a4c-2e03-4e0a-b9c5-bb6b51823999.6a86a0befb441ebc87fb0e036ebc5d5587d6ddb4.zh-cn.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "c28f4a4c-2e03-4e0a-b9c5-bb6b51823999.6a86a0befb441ebc87fb0e036ebc5d5587d6ddb4.zh-cn.xlf"
) | Out-Null

# ------------------------------------------------------------------
# de-de sheet: Latest Target File (F) / Latest Handback File (G)
# ------------------------------------------------------------------
$wsDeDe.Range("F2").Value = "23d4fd98-57d0-4775-a293-fe495abf9569.md"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/656fb5ebc286263f4ec67e3f3953099ed9a36ebf/e2e/23d4fd98-57d0-4775-a293-fe495abf9569.md",
    [Type]::Missing,
    [Type]::Missing,
    "23d4fd98-57d0-4775-a293-fe495abf9569.md"
) | Out-Null

$wsDeDe.Range("G2").Value = "23d4fd98-57d0-4775-a293-fe495abf9569.891c8a52f5279b6d8910c11a3544eaa27790e0f8.de-de.xlf"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/92afc32aeafae5d29614345f777e8bc392b3c0cf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/23d4fd98-57d0-4775-a293-fe495abf9569.891c8a52f5279b6d8910c11a3544eaa27790e0f8.de-de.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "23d4fd98-57d0-4775-a293-fe495abf9569.891c8a52f5279b6d8910c11a3544eaa27790e0f8.de-de.xlf"
) | Out-Null

$wsDeDe.Range("F3").Value = "c28f4a4c-2e03-4e0a-b9c5-bb6b51823999.md"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/656fb5ebc286263f4ec67e3f3953099ed9a36ebf/e2e/c28f4a4c-2e03-4e0a-b9c5-bb6b51823999.md",
    [Type]::Missing,
    [Type]::Missing,
    "c28f4a4c-2e03-4e0a-b9c5-bb6b51823999.md"
) | Out-Null

$wsDeDe.Range("G3").Value = "c28f4a4c-2e03-4e0a-b9c5-bb6b51823999.6a86a0befb441ebc87fb0e036ebc5d5587d6ddb4.de-de.xlf"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/92afc32aeafae5d29614345f777e8bc392b3c0cf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/c28f4a4c-2e03-4e0a-b9c5-bb6b51823999.6a86a0befb441ebc87fb0e036ebc5d5587d6ddb4.de-de.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "c28f4a4c-2e03-4e0a-b9c5-bb6b51823999.6a86a0befb441ebc87fb0e036ebc5d5587d6ddb4.de-de.xlf"
) | Out-Null

"Handback report generated."
